$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price/Volume columns so that numeric-looking
# strings (e.g. "229.74") are stored as text, matching the inlineStr cells in the source file.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.667.45'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '2.227.51'
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '0.643'
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").Value = '229.74'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = '63.19'
$ws.Range("E7").Value = '  +3.69%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.446'
$ws.Range("E9").Value = '  +5.11%  '
$ws.Range("D10").Value = '0.0964'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("D11").Value = '56.69'
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").Value = '26.59'
$ws.Range("E12").Value = '  +8.49%  '
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '2.566.92'
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").Value = '15.35'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '6.08'
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("D17").Value = '0.824'
$ws.Range("E17").Value = '  +1.47%  '
$ws.Range("D18").Value = '2.230.41'
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D19").Value = '43.589.41'
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("E20").Value = '  +3.98%  '
$ws.Range("D21").Value = '72.51'
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("E22").Value = '  -4.04%  '
$ws.Range("D23").Value = '248.52'
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -5.74%  '
$ws.Range("D26").Value = '3.39'
$ws.Range("E26").Value = '  +23.03%  '
$ws.Range("D27").Value = '2.29'
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("D29").Value = '170.47'
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("D30").Value = '20.77'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("E33").Value = '  +2.49%  '
$ws.Range("E34").Value = '  +5.54%  '
$ws.Range("D35").Value = '4.74'
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = '4.86'
$ws.Range("E36").Value = '  -4.09%  '
$ws.Range("D37").Value = '3.64'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").Value = '6.36'
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("D39").Value = '2.25'
$ws.Range("E39").Value = '  -6.03%  '
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").Value = '0.000218'
$ws.Range("E42").Value = '  -2.50%  '
$ws.Range("E43").Value = '  -7.54%  '
$ws.Range("D44").Value = '17.03'
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").Value = '96.72'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.35'
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '1.18'
$ws.Range("E47").Value = '  -2.92%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0938'
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("D49").Value = '2.33'
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.424.20'
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.76'
$ws.Range("E51").Value = '  +1.50%  '

# Reset style back to the default (Normal) so no stray number-format style is left
# behind on cells, matching the original (unstyled) cells in the source file.
$ws.Range("D2:E51").Style = "Normal"
